$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 526.9091
$ws.Range("I28").Value = 351.4074
$ws.Range("K28").Value = 351.4074
$ws.Range("M28").Value = 133.5926
$ws.Range("H108").Value = 34827.5
$ws.Range("J108").Value = 34827.5
$ws.Range("L108").Value = 34827.5
$ws.Range("N108").Value = -42507.5
$ws.Range("H112").Value = 1048.0646
$ws.Range("I112").Value = 494
$ws.Range("J112").Value = 1567.5
$ws.Range("K112").Value = 1482
$ws.Range("L112").Value = 4702.5
$ws.Range("M112").Value = -374
$ws.Range("N112").Value = -6918.5
$ws.Range("H124").Value = 36500
$ws.Range("J124").Value = 36500
$ws.Range("L124").Value = 36500
$ws.Range("N124").Value = -46320
$ws.Range("H129").Value = 718
$ws.Range("I129").Value = 539.4545000000001
$ws.Range("J129").Value = 1700
$ws.Range("K129").Value = 1618.3635
$ws.Range("L129").Value = 5100
$ws.Range("M129").Value = 3381.6365
$ws.Range("N129").Value = -15100
$ws.Range("H130").Value = 39172.855
$ws.Range("J130").Value = 39172.855
$ws.Range("L130").Value = 39172.855
$ws.Range("N130").Value = -49212.855
$ws.Range("H137").Value = 2253.1667
$ws.Range("I137").Value = 1169.7858
$ws.Range("K137").Value = 3509.3574
$ws.Range("M137").Value = -959.3574000000003
$ws.Range("H138").Value = 6316.5654
$ws.Range("I138").Value = 3966.2727
$ws.Range("J138").Value = 6635.7407
$ws.Range("K138").Value = 11898.8181
$ws.Range("L138").Value = 19907.2221
$ws.Range("M138").Value = -6758.8181
$ws.Range("N138").Value = -30187.2221

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1621.24
$ws.Range("I61").Value = 1566.5652
$ws.Range("K61").Value = 1566.5652
$ws.Range("M61").Value = -1354.5652
$ws.Range("H80").Value = 39888
$ws.Range("J80").Value = 39888
$ws.Range("L80").Value = 39888
$ws.Range("N80").Value = -41884
$ws.Range("H83").Value = 39888
$ws.Range("J83").Value = 39888
$ws.Range("L83").Value = 119664
$ws.Range("N83").Value = -129648
$ws.Range("H125").Value = 900000000
$ws.Range("J125").Value = 900000000
$ws.Range("L125").Value = 900000000
$ws.Range("N125").Value = -900009840
$ws.Range("H135").Value = 37183.727
$ws.Range("J135").Value = 37183.727
$ws.Range("L135").Value = 37183.727
$ws.Range("N135").Value = -47323.727
$ws.Range("H136").Value = 1621.24
$ws.Range("I136").Value = 1566.5652
$ws.Range("K136").Value = 4699.6956
$ws.Range("M136").Value = -2149.6956

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 68623.875
$ws.Range("I31").Value = 128025.25
$ws.Range("J31").Value = 9222.5
$ws.Range("K31").Value = 128025.25
$ws.Range("L31").Value = 9222.5
$ws.Range("M31").Value = -127730.25
$ws.Range("N31").Value = -9812.5
$ws.Range("H34").Value = 68623.875
$ws.Range("I34").Value = 128025.25
$ws.Range("J34").Value = 9222.5
$ws.Range("K34").Value = 128025.25
$ws.Range("L34").Value = 9222.5
$ws.Range("M34").Value = -127823.25
$ws.Range("N34").Value = -9626.5
$ws.Range("H50").Value = 7994
$ws.Range("J50").Value = 9312.5
$ws.Range("L50").Value = 9312.5
$ws.Range("N50").Value = -10562.5
$ws.Range("H60").Value = 12701.143
$ws.Range("J60").Value = 12701.143
$ws.Range("L60").Value = 12701.143
$ws.Range("N60").Value = -13723.143
$ws.Range("H68").Value = 15028
$ws.Range("J68").Value = 17492.666
$ws.Range("L68").Value = 17492.666
$ws.Range("N68").Value = -18990.666
$ws.Range("H71").Value = 15028
$ws.Range("J71").Value = 17492.666
$ws.Range("L71").Value = 52477.99800000001
$ws.Range("N71").Value = -59965.99800000001
$ws.Range("H127").Value = 41690
$ws.Range("J127").Value = 41690
$ws.Range("L127").Value = 41690
$ws.Range("N127").Value = -51610
$ws.Range("H141").Value = 63541.96
$ws.Range("J141").Value = 63541.96
$ws.Range("L141").Value = 63541.96
$ws.Range("N141").Value = -73901.95999999999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5550
$ws.Range("I133").Value = 5550
$ws.Range("K133").Value = 16650
$ws.Range("M133").Value = -11590
$ws.Range("H137").Value = 3025.889
$ws.Range("J137").Value = 4466.5
$ws.Range("L137").Value = 13399.5
$ws.Range("N137").Value = -23599.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1790
$ws.Range("I43").Value = 898.8889
$ws.Range("J43").Value = 9810
$ws.Range("K43").Value = 898.8889
$ws.Range("L43").Value = 9810
$ws.Range("M43").Value = -747.8889
$ws.Range("N43").Value = -10112
$ws.Range("H124").Value = 27461.111
$ws.Range("J124").Value = 27461.111
$ws.Range("L124").Value = 27461.111
$ws.Range("N124").Value = -37281.111
$ws.Range("H128").Value = 45650
$ws.Range("J128").Value = 45650
$ws.Range("L128").Value = 45650
$ws.Range("N128").Value = -55610
$ws.Range("H133").Value = 24745.166
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 24745.166
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 24745.166
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -34865.166
$ws.Range("H135").Value = 46880.91
$ws.Range("J135").Value = 46880.91
$ws.Range("L135").Value = 46880.91
$ws.Range("N135").Value = -57020.91

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -5376
$ws.Range("H63").Value = 32485
$ws.Range("J63").Value = 32485
$ws.Range("L63").Value = 32485
$ws.Range("N63").Value = -33983
$ws.Range("H66").Value = 32485
$ws.Range("J66").Value = 32485
$ws.Range("L66").Value = 97455
$ws.Range("N66").Value = -104943
$ws.Range("H125").Value = 42000
$ws.Range("J125").Value = 42000
$ws.Range("L125").Value = 42000
$ws.Range("N125").Value = -51840
$ws.Range("H134").Value = 61809.668
$ws.Range("J134").Value = 61809.668
$ws.Range("L134").Value = 61809.668
$ws.Range("N134").Value = -71949.66800000001
$ws.Range("H136").Value = 4503.273
$ws.Range("I136").Value = 1307.2858
$ws.Range("J136").Value = 22400.8
$ws.Range("K136").Value = 3921.8574
$ws.Range("L136").Value = 67202.39999999999
$ws.Range("M136").Value = -1371.8574
$ws.Range("N136").Value = -72302.39999999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 20127
$ws.Range("J109").Value = 20127
$ws.Range("L109").Value = 20127
$ws.Range("N109").Value = -22901
$ws.Range("H123").Value = 50356.855
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50356.855
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50356.855
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -60156.855
$ws.Range("H125").Value = 27357.5
$ws.Range("J125").Value = 27357.5
$ws.Range("L125").Value = 27357.5
$ws.Range("N125").Value = -37197.5
$ws.Range("H128").Value = 180397.86
$ws.Range("J128").Value = 180397.86
$ws.Range("L128").Value = 180397.86
$ws.Range("N128").Value = -190357.86
